$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 23:44"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 5908910
$ws.Range("C4").Value = 34764
$ws.Range("D4").Value = 3202187
$ws.Range("E4").Value = 2525771
$ws.Range("G4").Value = 348
$ws.Range("H4").Value = 180952

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 3622861
$ws.Range("C5").Value = 17078
$ws.Range("E5").Value = 797914
$ws.Range("G5").Value = 537
$ws.Range("H5").Value = 115309

# --- Row 8: Sudafrica ---
$ws.Range("B8").Value = 611450
$ws.Range("C8").Value = 1677
$ws.Range("D8").Value = 516494
$ws.Range("E8").Value = 81797
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = 13159

# --- Row 23: Alemania ---
$ws.Range("B23").Value = 236117
$ws.Range("C23").Value = 1628
$ws.Range("E23").Value = 17181

# --- Row 79: Costa de Marfil ---
$ws.Range("B79").Value = 17506
$ws.Range("C79").Value = 35
$ws.Range("D79").Value = 15633
$ws.Range("E79").Value = 1759
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 114

# --- Row 81: Bulgaria ---
$ws.Range("B81").Value = 15386
$ws.Range("C81").Value = 99
$ws.Range("D81").Value = 10497
$ws.Range("E81").Value = 4326
$ws.Range("G81").Value = 18
$ws.Range("H81").Value = 563

# --- Rows 123-125: Ruanda moves above Somalia/Mayotte with fresh data; Somalia and
# --- Mayotte shift down one row each, keeping their (former) numbers ---
$ws.Range("A123").Value = "Ruanda"
$ws.Range("B123").Value = 3306
$ws.Range("C123").Value = 217
$ws.Range("D123").Value = 1785
$ws.Range("E123").Value = 1507
$ws.Range("G123").Value = 2
$ws.Range("H123").Value = 14

$ws.Range("A124").Value = "Somalia"
$ws.Range("B124").Value = 3269
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 2396
$ws.Range("E124").Value = 780
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 93

$ws.Range("A125").Value = "Mayotte"
$ws.Range("B125").Value = 3237
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 2964
$ws.Range("E125").Value = 234
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 39

# --- Row 141: Yemen ---
$ws.Range("B141").Value = 1916
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 1090
$ws.Range("E141").Value = 271
$ws.Range("G141").Value = 2
$ws.Range("H141").Value = 555

# --- Row 156: Reunion ---
$ws.Range("B156").Value = 1244
$ws.Range("C156").Value = 35
$ws.Range("E156").Value = 546

# --- Row 162: Republica del Chad ---
$ws.Range("B162").Value = 987
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 41

# --- Rows 166-168: Belice moves above Crucero/San Marino with fresh data; Crucero and
# --- San Marino shift down one row each, keeping their (former) numbers ---
$ws.Range("A166").Value = "Belice"
$ws.Range("B166").Value = 713
$ws.Range("C166").Value = 27
$ws.Range("D166").Value = 45
$ws.Range("E166").Value = 658
$ws.Range("G166").Value = 4
$ws.Range("H166").Value = 10

$ws.Range("A167").Value = "Crucero"
$ws.Range("B167").Value = 712
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 651
$ws.Range("E167").Value = 48
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 13

$ws.Range("A168").Value = "San Marino"
$ws.Range("B168").Value = 704
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 657
$ws.Range("E168").Value = 5
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 42

# --- Rows 214-215: Montserrat and Islas Malvinas swap places ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
